$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 9 - SIAI blog post title/link update
$ws.Range("D9").Value = "SIAI의 Machine Learning 과목 기말 Term paper"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/then-prove-dnn-is-wrong/#utm_source=rss&utm_medium=rss&utm_campaign=then-prove-dnn-is-wrong"

# Row 20 - ai-creator tistory post title/link update
$ws.Range("D20").Value = "[책] [AI/MLOps] 머신러닝 파워드 애플리케이션(Building Machine Learning Powered Application)"
$ws.Range("E20").Value = "https://ai-creator.tistory.com/617"

# Row 36 - dmqm seminar title/link update
$ws.Range("D36").Value = "Anomaly Detection for Time Series with Autoencoder"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/358"

# Row 37 - dsba seminar title update (link unchanged)
$ws.Range("D37").Value = "dsba_seminar"
